$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: PRIORITY=High, STATUS=In Progress, START DATE=7/1/2014, % COMPLETE=25%
$ws.Range("D10").Value = "High"
$ws.Range("E10").Value = "In Progress"
$ws.Range("F10").Value = (Get-Date -Year 2014 -Month 7 -Day 1)
$ws.Range("H10").Value = 0.25

# Row 11: PRIORITY=High, STATUS=In Progress, START DATE=7/1/2014, % COMPLETE=25%
$ws.Range("D11").Value = "High"
$ws.Range("E11").Value = "In Progress"
$ws.Range("F11").Value = (Get-Date -Year 2014 -Month 7 -Day 1)
$ws.Range("H11").Value = 0.25

# Row 12: PRIORITY=Normal
$ws.Range("D12").Value = "Normal"

# Row 13: PRIORITY=Normal
$ws.Range("D13").Value = "Normal"

# Row 14: PRIORITY=High
$ws.Range("D14").Value = "High"

# Update selection to H11
$ws.Range("H11").Select()
